# Auto-generated edit script: refresh market-price-derived columns (H-N)
# on the per-job Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR),
# mirroring the scheduled price-sync job's output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2305.5715
$ws.Range("I2").Value = 429.75
$ws.Range("K2").Value = 429.75
$ws.Range("M2").Value = -316.75

$ws.Range("H4").Value = 470.75
$ws.Range("I4").Value = 461.33334
$ws.Range("J4").Value = 499
$ws.Range("K4").Value = 461.33334
$ws.Range("L4").Value = 499
$ws.Range("M4").Value = -347.33334
$ws.Range("N4").Value = -727

$ws.Range("H112").Value = 3061.7932
$ws.Range("J112").Value = 3130.4443
$ws.Range("L112").Value = 9391.332900000001
$ws.Range("N112").Value = -11607.3329

$ws.Range("H132").Value = 502406.25
$ws.Range("I132").Value = 2565.5151
$ws.Range("J132").Value = 2858798.2
$ws.Range("K132").Value = 7696.5453
$ws.Range("L132").Value = 8576394.600000001
$ws.Range("M132").Value = -5166.5453
$ws.Range("N132").Value = -8581454.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18795.848
$ws.Range("I32").Value = 6631.6045
$ws.Range("K32").Value = 6631.6045
$ws.Range("M32").Value = -6344.6045

$ws.Range("H122").Value = 3110.4688
$ws.Range("I122").Value = 3045.44
$ws.Range("K122").Value = 9136.32
$ws.Range("M122").Value = -6686.32

$ws.Range("H132").Value = 2562.4614
$ws.Range("I132").Value = 1109.5
$ws.Range("K132").Value = 3328.5
$ws.Range("M132").Value = -798.5

$ws.Range("H139").Value = 72530.57000000001
$ws.Range("J139").Value = 72530.57000000001
$ws.Range("L139").Value = 72530.57000000001
$ws.Range("N139").Value = -82810.57000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 90000
$ws.Range("J104").Value = 90000
$ws.Range("L104").Value = 90000
$ws.Range("N104").Value = -96988

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 869.65216
$ws.Range("I107").Value = 860.2143
$ws.Range("K107").Value = 860.2143
$ws.Range("M107").Value = 1059.7857

$ws.Range("H122").Value = 3106.875
$ws.Range("I122").Value = 3272.7856
$ws.Range("K122").Value = 9818.356800000001
$ws.Range("M122").Value = -7368.356800000001

$ws.Range("H134").Value = 2173.75
$ws.Range("I134").Value = 2312.8572
$ws.Range("K134").Value = 6938.571599999999
$ws.Range("M134").Value = -4403.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 57.545456
$ws.Range("I12").Value = 61.75
$ws.Range("K12").Value = 185.25
$ws.Range("M12").Value = -12.25

$ws.Range("H39").Value = 3345.923
$ws.Range("I39").Value = 2999
$ws.Range("K39").Value = 8997
$ws.Range("M39").Value = -8703

$ws.Range("H55").Value = 3377.4546
$ws.Range("J55").Value = 5070.5713
$ws.Range("L55").Value = 15211.7139
$ws.Range("N55").Value = -15565.7139

$ws.Range("H69").Value = 9991.9
$ws.Range("I69").Value = 1949
$ws.Range("J69").Value = 12002.625
$ws.Range("K69").Value = 5847
$ws.Range("L69").Value = 36007.875
$ws.Range("M69").Value = -5036
$ws.Range("N69").Value = -37629.875

$ws.Range("H72").Value = 9991.9
$ws.Range("I72").Value = 1949
$ws.Range("J72").Value = 12002.625
$ws.Range("K72").Value = 17541
$ws.Range("L72").Value = 108023.625
$ws.Range("M72").Value = -13485
$ws.Range("N72").Value = -116135.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3112.0588
$ws.Range("I102").Value = 2883.4167
$ws.Range("J102").Value = 3660.8
$ws.Range("K102").Value = 2883.4167
$ws.Range("L102").Value = 3660.8
$ws.Range("M102").Value = -1261.4167
$ws.Range("N102").Value = -6904.8

$ws.Range("H126").Value = 3206.1333
$ws.Range("I126").Value = 3024
$ws.Range("K126").Value = 9072
$ws.Range("M126").Value = -6602

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6179.3076
$ws.Range("I7").Value = 5036.7915
$ws.Range("K7").Value = 5036.7915
$ws.Range("M7").Value = -4924.7915

$ws.Range("H40").Value = 5331.2856
$ws.Range("I40").Value = 5775.6665
$ws.Range("J40").Value = 4998
$ws.Range("K40").Value = 5775.6665
$ws.Range("L40").Value = 4998
$ws.Range("M40").Value = -5639.6665
$ws.Range("N40").Value = -5270

$ws.Range("H61").Value = 1256.7059
$ws.Range("J61").Value = 791
$ws.Range("L61").Value = 791
$ws.Range("N61").Value = -1195

$ws.Range("H108").Value = 89981.5
$ws.Range("J108").Value = 89981.5
$ws.Range("L108").Value = 89981.5
$ws.Range("N108").Value = -97661.5

$ws.Range("H113").Value = 1256.7059
$ws.Range("J113").Value = 791
$ws.Range("L113").Value = 791
$ws.Range("N113").Value = -5131

$ws.Range("H119").Value = 97514.5
$ws.Range("J119").Value = 97514.5
$ws.Range("L119").Value = 97514.5
$ws.Range("N119").Value = -107190.5

$ws.Range("H120").Value = 110000
$ws.Range("J120").Value = 110000
$ws.Range("L120").Value = 110000
$ws.Range("N120").Value = -119676

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 3072.6667
$ws.Range("I122").Value = 2660.125
$ws.Range("J122").Value = 3672.7273
$ws.Range("K122").Value = 7980.375
$ws.Range("L122").Value = 11018.1819
$ws.Range("M122").Value = -5530.375
$ws.Range("N122").Value = -15918.1819

$ws.Range("H123").Value = 68250
$ws.Range("J123").Value = 68250
$ws.Range("L123").Value = 68250
$ws.Range("N123").Value = -78050

$ws.Range("H126").Value = 6179.3076
$ws.Range("I126").Value = 5036.7915
$ws.Range("K126").Value = 15110.3745
$ws.Range("M126").Value = -12640.3745

$ws.Range("H136").Value = 4318.1665
$ws.Range("I136").Value = 4167.4287
$ws.Range("K136").Value = 12502.2861
$ws.Range("M136").Value = -9952.286100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 59997.5
$ws.Range("J108").Value = 59997.5
$ws.Range("L108").Value = 59997.5
$ws.Range("N108").Value = -67677.5

$ws.Range("H109").Value = 76075.39999999999
$ws.Range("J109").Value = 76075.39999999999
$ws.Range("L109").Value = 76075.39999999999
$ws.Range("N109").Value = -78849.39999999999

$ws.Range("H126").Value = 3796.6191
$ws.Range("I126").Value = 2807.7058
$ws.Range("J126").Value = 7999.5
$ws.Range("K126").Value = 8423.117400000001
$ws.Range("L126").Value = 23998.5
$ws.Range("M126").Value = -5953.117400000001
$ws.Range("N126").Value = -28938.5
